$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (A7 = 11 / Q6)
$ws.Range("B7").Value = -0.09360001116797961
$ws.Range("C7").Value = 0.6773215120173049
$ws.Range("D7").Value = 0.8387150314969378
$ws.Range("E7").Value = 0.9158138629093457
$ws.Range("F7").Value = 0.9232471286017232
$ws.Range("G7").Value = 38

# Row 8 (A8 = 12 / Q7)
$ws.Range("B8").Value = -0.07158090847683796
$ws.Range("C8").Value = 0.5716789624051219
$ws.Range("D8").Value = 0.6667522629364007
$ws.Range("E8").Value = 0.8165489960415117
$ws.Range("F8").Value = 0.8246253720071778
$ws.Range("G8").Value = 37

# Row 9 (A9 = 13 / Q8)
$ws.Range("B9").Value = -0.347548237505526
$ws.Range("C9").Value = 0.4489756673986917
$ws.Range("D9").Value = 0.2984715574757459
$ws.Range("E9").Value = 0.5463255050569632
$ws.Range("F9").Value = 0.4324736439581862
$ws.Range("G9").Value = 20

# Row 10 (A10 = 14 / Q9)
$ws.Range("B10").Value = -0.2871006105435207
$ws.Range("C10").Value = 0.4966077703424752
$ws.Range("D10").Value = 0.3494510458411312
$ws.Range("E10").Value = 0.5911438453042805
$ws.Range("F10").Value = 0.5378441308026808
$ws.Range("G10").Value = 13

# Row 11 (A11 = 15 / last row)
$ws.Range("B11").Value = -0.3134163365220089
$ws.Range("C11").Value = 0.3880630270395826
$ws.Range("D11").Value = 0.2215234913900424
$ws.Range("E11").Value = 0.4706628213382085
$ws.Range("F11").Value = 0.3925775264058765
$ws.Range("G11").Value = 5
